# Apply cell value updates to multiple sheets per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 157.5
$ws.Range("I11").Value = 157.5
$ws.Range("K11").Value = 157.5
$ws.Range("M11").Value = -17.5
$ws.Range("H17").Value = 874.7143
$ws.Range("J17").Value = 874.7143
$ws.Range("L17").Value = 2624.1429
$ws.Range("N17").Value = -2960.1429
$ws.Range("H106").Value = 8914.177
$ws.Range("I106").Value = 10102.5
$ws.Range("J106").Value = 3368.6667
$ws.Range("K106").Value = 10102.5
$ws.Range("L106").Value = 3368.6667
$ws.Range("M106").Value = -9471.5
$ws.Range("N106").Value = -4630.6667
$ws.Range("H127").Value = 1840
$ws.Range("I127").Value = 350
$ws.Range("J127").Value = 3330
$ws.Range("K127").Value = 1050
$ws.Range("L127").Value = 9990
$ws.Range("M127").Value = 3910
$ws.Range("N127").Value = -19910
$ws.Range("H132").Value = 8551813
$ws.Range("I132").Value = 9264092
$ws.Range("J132").Value = 4470.6665
$ws.Range("K132").Value = 27792276
$ws.Range("L132").Value = 13411.9995
$ws.Range("M132").Value = -27789746
$ws.Range("N132").Value = -18471.9995
$ws.Range("H137").Value = 1128.5667
$ws.Range("I137").Value = 896.5454999999999
$ws.Range("J137").Value = 1766.625
$ws.Range("K137").Value = 2689.6365
$ws.Range("L137").Value = 5299.875
$ws.Range("M137").Value = -139.6364999999996
$ws.Range("N137").Value = -10399.875
$ws.Range("H141").Value = 1675
$ws.Range("I141").Value = 350
$ws.Range("K141").Value = 1050
$ws.Range("M141").Value = 4130

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4367.8975
$ws.Range("I32").Value = 3953.7222
$ws.Range("K32").Value = 3953.7222
$ws.Range("M32").Value = -3666.7222
$ws.Range("H61").Value = 863.0323
$ws.Range("I61").Value = 660.5833
$ws.Range("J61").Value = 1557.1428
$ws.Range("K61").Value = 660.5833
$ws.Range("L61").Value = 1557.1428
$ws.Range("M61").Value = -448.5833
$ws.Range("N61").Value = -1981.1428
$ws.Range("H74").Value = 1071.3438
$ws.Range("I74").Value = 567.73914
$ws.Range("J74").Value = 2358.3333
$ws.Range("K74").Value = 567.73914
$ws.Range("L74").Value = 2358.3333
$ws.Range("M74").Value = 306.26086
$ws.Range("N74").Value = -4106.3333
$ws.Range("H77").Value = 1071.3438
$ws.Range("I77").Value = 567.73914
$ws.Range("J77").Value = 2358.3333
$ws.Range("K77").Value = 2838.6957
$ws.Range("L77").Value = 11791.6665
$ws.Range("M77").Value = 1529.3043
$ws.Range("N77").Value = -20527.6665
$ws.Range("H136").Value = 863.0323
$ws.Range("I136").Value = 660.5833
$ws.Range("J136").Value = 1557.1428
$ws.Range("K136").Value = 1981.7499
$ws.Range("L136").Value = 4671.428400000001
$ws.Range("M136").Value = 568.2501
$ws.Range("N136").Value = -9771.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 397.25
$ws.Range("I64").Value = 411.77777
$ws.Range("J64").Value = 378.57144
$ws.Range("K64").Value = 411.77777
$ws.Range("L64").Value = 378.57144
$ws.Range("M64").Value = -186.77777
$ws.Range("N64").Value = -828.5714399999999
$ws.Range("H67").Value = 397.25
$ws.Range("I67").Value = 411.77777
$ws.Range("J67").Value = 378.57144
$ws.Range("K67").Value = 411.77777
$ws.Range("L67").Value = 378.57144
$ws.Range("M67").Value = 368.22223
$ws.Range("N67").Value = -1938.57144
$ws.Range("H105").Value = 125002590
$ws.Range("I105").Value = 125002590
$ws.Range("K105").Value = 125002590
$ws.Range("M105").Value = -125000843

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1588.7037
$ws.Range("I31").Value = 1248.2632
$ws.Range("K31").Value = 1248.2632
$ws.Range("M31").Value = -953.2632000000001
$ws.Range("H34").Value = 1588.7037
$ws.Range("I34").Value = 1248.2632
$ws.Range("K34").Value = 1248.2632
$ws.Range("M34").Value = -1046.2632
$ws.Range("H58").Value = 1123.9445
$ws.Range("I58").Value = 810.1667
$ws.Range("K58").Value = 810.1667
$ws.Range("M58").Value = -607.1667
$ws.Range("H99").Value = 1846
$ws.Range("I99").Value = 1846
$ws.Range("K99").Value = 1846
$ws.Range("M99").Value = -348
$ws.Range("H126").Value = 1846
$ws.Range("I126").Value = 1846
$ws.Range("K126").Value = 5538
$ws.Range("M126").Value = -3068
$ws.Range("H132").Value = 4907.6665
$ws.Range("I132").Value = 6262.5713
$ws.Range("J132").Value = 2536.5833
$ws.Range("K132").Value = 18787.7139
$ws.Range("L132").Value = 7609.749899999999
$ws.Range("M132").Value = -16257.7139
$ws.Range("N132").Value = -12669.7499
$ws.Range("H136").Value = 1123.9445
$ws.Range("I136").Value = 810.1667
$ws.Range("K136").Value = 2430.5001
$ws.Range("M136").Value = 119.4998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2795.4443
$ws.Range("I136").Value = 1676.6666
$ws.Range("K136").Value = 5029.9998
$ws.Range("M136").Value = 70.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 45003410
$ws.Range("I70").Value = 35717916
$ws.Range("J70").Value = 66669570
$ws.Range("K70").Value = 35717916
$ws.Range("L70").Value = 66669570
$ws.Range("M70").Value = -35717646
$ws.Range("N70").Value = -66670110
$ws.Range("H73").Value = 45003410
$ws.Range("I73").Value = 35717916
$ws.Range("J73").Value = 66669570
$ws.Range("K73").Value = 35717916
$ws.Range("L73").Value = 66669570
$ws.Range("M73").Value = -35716980
$ws.Range("N73").Value = -66671442
$ws.Range("H80").Value = 5138.2856
$ws.Range("I80").Value = 4667
$ws.Range("J80").Value = 5766.6665
$ws.Range("K80").Value = 4667
$ws.Range("L80").Value = 5766.6665
$ws.Range("M80").Value = -3669
$ws.Range("N80").Value = -7762.6665
$ws.Range("H83").Value = 5138.2856
$ws.Range("I83").Value = 4667
$ws.Range("J83").Value = 5766.6665
$ws.Range("K83").Value = 23335
$ws.Range("L83").Value = 28833.3325
$ws.Range("M83").Value = -18343
$ws.Range("N83").Value = -38817.3325
$ws.Range("H132").Value = 1959.5625
$ws.Range("I132").Value = 1756.32
$ws.Range("J132").Value = 2685.4285
$ws.Range("K132").Value = 5268.96
$ws.Range("L132").Value = 8056.2855
$ws.Range("M132").Value = -2738.96
$ws.Range("N132").Value = -13116.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1634.826
$ws.Range("I132").Value = 1225.4688
$ws.Range("J132").Value = 2570.5
$ws.Range("K132").Value = 3676.4064
$ws.Range("L132").Value = 7711.5
$ws.Range("M132").Value = -1146.4064
$ws.Range("N132").Value = -12771.5
